$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1122.1364
$ws.Range("J17").Value = 1122.1364
$ws.Range("L17").Value = 3366.4092
$ws.Range("N17").Value = -3702.4092

$ws.Range("H62").Value = 2892.1667
$ws.Range("I62").Value = 2299.25
$ws.Range("K62").Value = 2299.25
$ws.Range("M62").Value = -1675.25

$ws.Range("H65").Value = 2892.1667
$ws.Range("I65").Value = 2299.25
$ws.Range("K65").Value = 11496.25
$ws.Range("M65").Value = -8376.25

$ws.Range("H116").Value = 20002760
$ws.Range("I116").Value = 50002024
$ws.Range("J116").Value = 3250
$ws.Range("K116").Value = 50002024
$ws.Range("L116").Value = 3250
$ws.Range("M116").Value = -49998582
$ws.Range("N116").Value = -10134

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1233.6897
$ws.Range("I61").Value = 1082.52
$ws.Range("J61").Value = 2178.5
$ws.Range("K61").Value = 1082.52
$ws.Range("L61").Value = 2178.5
$ws.Range("M61").Value = -870.52
$ws.Range("N61").Value = -2602.5

$ws.Range("H74").Value = 1096.4073
$ws.Range("I74").Value = 842.8421
$ws.Range("J74").Value = 1698.625
$ws.Range("K74").Value = 842.8421
$ws.Range("L74").Value = 1698.625
$ws.Range("M74").Value = 31.15790000000004
$ws.Range("N74").Value = -3446.625

$ws.Range("H77").Value = 1096.4073
$ws.Range("I77").Value = 842.8421
$ws.Range("J77").Value = 1698.625
$ws.Range("K77").Value = 4214.2105
$ws.Range("L77").Value = 8493.125
$ws.Range("M77").Value = 153.7894999999999
$ws.Range("N77").Value = -17229.125

$ws.Range("H96").Value = 120488.8
$ws.Range("J96").Value = 120488.8
$ws.Range("L96").Value = 120488.8
$ws.Range("N96").Value = -125980.8

$ws.Range("H132").Value = 2989.8
$ws.Range("I132").Value = 1816.6666
$ws.Range("J132").Value = 4749.5
$ws.Range("K132").Value = 5449.9998
$ws.Range("L132").Value = 14248.5
$ws.Range("M132").Value = -2919.9998
$ws.Range("N132").Value = -19308.5

$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

$ws.Range("H136").Value = 1233.6897
$ws.Range("I136").Value = 1082.52
$ws.Range("J136").Value = 2178.5
$ws.Range("K136").Value = 3247.56
$ws.Range("L136").Value = 6535.5
$ws.Range("M136").Value = -697.5599999999999
$ws.Range("N136").Value = -11635.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2292.6667
$ws.Range("I134").Value = 2128.5715
$ws.Range("J134").Value = 2522.4
$ws.Range("K134").Value = 6385.7145
$ws.Range("L134").Value = 7567.200000000001
$ws.Range("M134").Value = -3850.7145
$ws.Range("N134").Value = -12637.2

$ws.Range("H135").Value = 110080
$ws.Range("J135").Value = 110080
$ws.Range("L135").Value = 110080
$ws.Range("N135").Value = -120220

$ws.Range("H137").Value = 59487.5
$ws.Range("J137").Value = 59487.5
$ws.Range("L137").Value = 59487.5
$ws.Range("N137").Value = -69687.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 10316.5
$ws.Range("I69").Value = 9379.799999999999
$ws.Range("K69").Value = 9379.799999999999
$ws.Range("M69").Value = -8630.799999999999

$ws.Range("H72").Value = 10316.5
$ws.Range("I72").Value = 9379.799999999999
$ws.Range("K72").Value = 28139.4
$ws.Range("M72").Value = -24395.4

$ws.Range("H132").Value = 2225.8
$ws.Range("I132").Value = 1743.3334
$ws.Range("J132").Value = 2949.5
$ws.Range("K132").Value = 5230.0002
$ws.Range("L132").Value = 8848.5
$ws.Range("M132").Value = -2700.0002
$ws.Range("N132").Value = -13908.5

$ws.Range("H134").Value = 2836
$ws.Range("I134").Value = 2644.5715
$ws.Range("J134").Value = 3506
$ws.Range("K134").Value = 7933.7145
$ws.Range("L134").Value = 10518
$ws.Range("M134").Value = -5398.7145
$ws.Range("N134").Value = -15588

$ws.Range("H135").Value = 32868.75
$ws.Range("J135").Value = 32868.75
$ws.Range("L135").Value = 32868.75
$ws.Range("N135").Value = -43008.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1000644.1
$ws.Range("J113").Value = 693
$ws.Range("L113").Value = 2079
$ws.Range("N113").Value = -6419

$ws.Range("H141").Value = 4847.25
$ws.Range("I141").Value = 3362.9412
$ws.Range("J141").Value = 8452
$ws.Range("K141").Value = 10088.8236
$ws.Range("L141").Value = 25356
$ws.Range("M141").Value = -4908.8236
$ws.Range("N141").Value = -35716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7853.091
$ws.Range("I70").Value = 7625
$ws.Range("J70").Value = 7983.4287
$ws.Range("K70").Value = 7625
$ws.Range("L70").Value = 7983.4287
$ws.Range("M70").Value = -7355
$ws.Range("N70").Value = -8523.4287

$ws.Range("H73").Value = 7853.091
$ws.Range("I73").Value = 7625
$ws.Range("J73").Value = 7983.4287
$ws.Range("K73").Value = 7625
$ws.Range("L73").Value = 7983.4287
$ws.Range("M73").Value = -6689
$ws.Range("N73").Value = -9855.4287

$ws.Range("H92").Value = 32475
$ws.Range("J92").Value = 32475
$ws.Range("L92").Value = 32475
$ws.Range("N92").Value = -36219

$ws.Range("H95").Value = 1683531.4
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 1683531.4
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 1683531.4
$ws.Range("M95").ClearContents()
$ws.Range("N95").Value = -1689023.4

$ws.Range("H122").Value = 3927.5454
$ws.Range("I122").Value = 4274.375
$ws.Range("J122").Value = 3002.6667
$ws.Range("K122").Value = 12823.125
$ws.Range("L122").Value = 9008.000100000001
$ws.Range("M122").Value = -10373.125
$ws.Range("N122").Value = -13908.0001

$ws.Range("H132").Value = 2751.7896
$ws.Range("I132").Value = 2049.7273
$ws.Range("K132").Value = 6149.1819
$ws.Range("M132").Value = -3619.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1528.5714
$ws.Range("I93").Value = 1600
$ws.Range("J93").Value = 1475
$ws.Range("K93").Value = 1600
$ws.Range("L93").Value = 1475
$ws.Range("M93").Value = -352
$ws.Range("N93").Value = -3971

$ws.Range("H94").Value = 16565
$ws.Range("J94").Value = 16565
$ws.Range("L94").Value = 16565
$ws.Range("N94").Value = -17917

$ws.Range("H100").Value = 6126.8667
$ws.Range("I100").Value = 8712.875
$ws.Range("J100").Value = 3171.4285
$ws.Range("K100").Value = 8712.875
$ws.Range("L100").Value = 3171.4285
$ws.Range("M100").Value = -8171.875
$ws.Range("N100").Value = -4253.4285

$ws.Range("H136").Value = 3390.3462
$ws.Range("I136").Value = 2744.6843
$ws.Range("J136").Value = 5142.857
$ws.Range("K136").Value = 8234.052899999999
$ws.Range("L136").Value = 15428.571
$ws.Range("M136").Value = -5684.052899999999
$ws.Range("N136").Value = -20528.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2504.7144
$ws.Range("I132").Value = 1823.2858
$ws.Range("J132").Value = 3186.1428
$ws.Range("K132").Value = 5469.857400000001
$ws.Range("L132").Value = 9558.428400000001
$ws.Range("M132").Value = -2939.857400000001
$ws.Range("N132").Value = -14618.4284
